$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(28, 9).Value = "aa"
$ws.Cells.Item(28, 10).Value = "Agree/Accept"
$ws.Cells.Item(29, 9).Value = "%"
$ws.Cells.Item(29, 10).Value = "Uninterpretable"
$ws.Cells.Item(39, 9).Value = "aa"
$ws.Cells.Item(39, 10).Value = "Agree/Accept"
$ws.Cells.Item(44, 9).Value = "sd"
$ws.Cells.Item(44, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(46, 9).Value = "b"
$ws.Cells.Item(46, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(56, 9).Value = "sd"
$ws.Cells.Item(56, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(64, 9).Value = "sd"
$ws.Cells.Item(64, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(70, 9).Value = "sd"
$ws.Cells.Item(70, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(72, 9).Value = "sd"
$ws.Cells.Item(72, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(83, 9).Value = "b"
$ws.Cells.Item(83, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(112, 9).Value = "%"
$ws.Cells.Item(112, 10).Value = "Uninterpretable"
$ws.Cells.Item(113, 9).Value = "sv"
$ws.Cells.Item(113, 10).Value = "Statement-opinion"
$ws.Cells.Item(114, 9).Value = "sd"
$ws.Cells.Item(114, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(123, 9).Value = "sd"
$ws.Cells.Item(123, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(124, 9).Value = "sv"
$ws.Cells.Item(124, 10).Value = "Statement-opinion"
$ws.Cells.Item(135, 9).Value = "sv"
$ws.Cells.Item(135, 10).Value = "Statement-opinion"
$ws.Cells.Item(144, 9).Value = "sd"
$ws.Cells.Item(144, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(152, 9).Value = "sd"
$ws.Cells.Item(152, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(156, 9).Value = "b"
$ws.Cells.Item(156, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(163, 9).Value = "sd"
$ws.Cells.Item(163, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(164, 9).Value = "sd"
$ws.Cells.Item(164, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(168, 9).Value = "sv"
$ws.Cells.Item(168, 10).Value = "Statement-opinion"
$ws.Cells.Item(169, 9).Value = "sv"
$ws.Cells.Item(169, 10).Value = "Statement-opinion"
$ws.Cells.Item(183, 9).Value = "sv"
$ws.Cells.Item(183, 10).Value = "Statement-opinion"
$ws.Cells.Item(191, 9).Value = "sd"
$ws.Cells.Item(191, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(207, 9).Value = "sv"
$ws.Cells.Item(207, 10).Value = "Statement-opinion"
$ws.Cells.Item(222, 9).Value = "qy"
$ws.Cells.Item(222, 10).Value = "Yes-No-Question"
$ws.Cells.Item(225, 9).Value = "ba"
$ws.Cells.Item(225, 10).Value = "Appreciation"
$ws.Cells.Item(227, 9).Value = "sv"
$ws.Cells.Item(227, 10).Value = "Statement-opinion"
$ws.Cells.Item(230, 9).Value = "sd"
$ws.Cells.Item(230, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(238, 9).Value = "sd"
$ws.Cells.Item(238, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(239, 9).Value = "sv"
$ws.Cells.Item(239, 10).Value = "Statement-opinion"
$ws.Cells.Item(250, 9).Value = "b"
$ws.Cells.Item(250, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(254, 9).Value = "%"
$ws.Cells.Item(254, 10).Value = "Uninterpretable"
$ws.Cells.Item(255, 9).Value = "aa"
$ws.Cells.Item(255, 10).Value = "Agree/Accept"
$ws.Cells.Item(277, 9).Value = "sd"
$ws.Cells.Item(277, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(303, 9).Value = "sd"
$ws.Cells.Item(303, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(309, 9).Value = "%"
$ws.Cells.Item(309, 10).Value = "Uninterpretable"
$ws.Cells.Item(317, 9).Value = "%"
$ws.Cells.Item(317, 10).Value = "Uninterpretable"
$ws.Cells.Item(320, 9).Value = "sv"
$ws.Cells.Item(320, 10).Value = "Statement-opinion"
$ws.Cells.Item(339, 9).Value = "b"
$ws.Cells.Item(339, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(343, 9).Value = "sv"
$ws.Cells.Item(343, 10).Value = "Statement-opinion"
$ws.Cells.Item(345, 9).Value = "aa"
$ws.Cells.Item(345, 10).Value = "Agree/Accept"
$ws.Cells.Item(356, 9).Value = "ba"
$ws.Cells.Item(356, 10).Value = "Appreciation"
$ws.Cells.Item(367, 9).Value = "sv"
$ws.Cells.Item(367, 10).Value = "Statement-opinion"
$ws.Cells.Item(370, 9).Value = "b"
$ws.Cells.Item(370, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(378, 9).Value = "sv"
$ws.Cells.Item(378, 10).Value = "Statement-opinion"
$ws.Cells.Item(384, 9).Value = "sv"
$ws.Cells.Item(384, 10).Value = "Statement-opinion"
$ws.Cells.Item(387, 9).Value = "sd"
$ws.Cells.Item(387, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(396, 9).Value = "aa"
$ws.Cells.Item(396, 10).Value = "Agree/Accept"
$ws.Cells.Item(398, 9).Value = "aa"
$ws.Cells.Item(398, 10).Value = "Agree/Accept"
$ws.Cells.Item(403, 9).Value = "sd"
$ws.Cells.Item(403, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(415, 9).Value = "%"
$ws.Cells.Item(415, 10).Value = "Uninterpretable"
$ws.Cells.Item(417, 9).Value = "sd"
$ws.Cells.Item(417, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(425, 9).Value = "aa"
$ws.Cells.Item(425, 10).Value = "Agree/Accept"
$ws.Cells.Item(430, 9).Value = "sv"
$ws.Cells.Item(430, 10).Value = "Statement-opinion"
$ws.Cells.Item(432, 9).Value = "sd"
$ws.Cells.Item(432, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(433, 9).Value = "sd"
$ws.Cells.Item(433, 10).Value = "Statement-non-opinion"
